# Apply the cryptos-list refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain a text value even when it looks numeric
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '63.584.09'
$ws.Range("E2").Value = '  -1.14%  '
Set-TextValue $ws.Range("D3") '3.421.53'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("E4").Value = '  +0.12%  '
Set-TextValue $ws.Range("D5") '578.41'
$ws.Range("E5").Value = '  -2.17%  '
Set-TextValue $ws.Range("D6") '128.94'
$ws.Range("E6").Value = '  -3.90%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -1.54%  '
Set-TextValue $ws.Range("D9") '7.55'
$ws.Range("E9").Value = '  +2.77%  '
$ws.Range("E10").Value = '  -1.30%  '
Set-TextValue $ws.Range("D11") '0.381'
$ws.Range("E11").Value = '  -1.39%  '
Set-TextValue $ws.Range("D12") '4.007.65'
$ws.Range("E12").Value = '  -2.24%  '
$ws.Range("E13").Value = '  -0.36%  '
Set-TextValue $ws.Range("D14") '0.0000175'
$ws.Range("E14").Value = '  -3.11%  '
Set-TextValue $ws.Range("D15") '3.424.14'
$ws.Range("E15").Value = '  -2.23%  '
Set-TextValue $ws.Range("D16") '63.654.35'
$ws.Range("E16").Value = '  -1.16%  '
Set-TextValue $ws.Range("D17") '25.01'
$ws.Range("E17").Value = '  -2.61%  '
Set-TextValue $ws.Range("D18") '9.82'
$ws.Range("E18").Value = '  -0.51%  '
Set-TextValue $ws.Range("D19") '5.62'
$ws.Range("E19").Value = '  -2.29%  '
Set-TextValue $ws.Range("D20") '13.26'
$ws.Range("E20").Value = '  -2.00%  '
Set-TextValue $ws.Range("D21") '384.83'
$ws.Range("E21").Value = '  -2.17%  '
Set-TextValue $ws.Range("D22") '0.562'
$ws.Range("E22").Value = '  -1.98%  '
Set-TextValue $ws.Range("D23") '3.562.94'
$ws.Range("E23").Value = '  -2.18%  '
Set-TextValue $ws.Range("D24") '73.92'
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("E26").Value = '  -5.48%  '
Set-TextValue $ws.Range("D27") '0.999'
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("E28").Value = '  -3.37%  '
Set-TextValue $ws.Range("D29") '7.01'
$ws.Range("E29").Value = '  -5.08%  '
Set-TextValue $ws.Range("D30") '7.88'
$ws.Range("E30").Value = '  -4.31%  '
$ws.Range("E31").Value = '  -0.47%  '
Set-TextValue $ws.Range("D32") '1.40'
$ws.Range("E32").Value = '  -4.86%  '
Set-TextValue $ws.Range("D33") '3.457.28'
$ws.Range("E33").Value = '  -1.90%  '
Set-TextValue $ws.Range("D35") '22.66'
$ws.Range("E35").Value = '  -3.51%  '
Set-TextValue $ws.Range("D36") '5.16'
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("E37").Value = '  -2.20%  '
Set-TextValue $ws.Range("D38") '164.12'
$ws.Range("E38").Value = '  -1.97%  '
$ws.Range("E39").Value = '  -2.78%  '
Set-TextValue $ws.Range("D40") '0.0766'
$ws.Range("E40").Value = '  -1.86%  '
Set-TextValue $ws.Range("D41") '0.782'
$ws.Range("E41").Value = '  -3.20%  '
$ws.Range("E42").Value = '  -0.04%  '
Set-TextValue $ws.Range("D43") '41.34'
$ws.Range("E43").Value = '  -1.06%  '
Set-TextValue $ws.Range("D44") '4.30'
$ws.Range("E44").Value = '  -2.16%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D45") '23.40'
$ws.Range("E45").Value = '  -6.52%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D46") '1.59'
$ws.Range("E46").Value = '  -4.28%  '
$ws.Range("E47").Value = '  -6.20%  '
Set-TextValue $ws.Range("D48") '6.69'
$ws.Range("E48").Value = '  -0.97%  '
Set-TextValue $ws.Range("D49") '0.881'
$ws.Range("E49").Value = '  -1.58%  '
Set-TextValue $ws.Range("D50") '2.271.67'
$ws.Range("E50").Value = '  -4.52%  '
$ws.Range("E51").Value = '  -2.58%  '
